# Update PLC data 2025-10-13 14:14:03
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 177621
$ws.Range("C4").Value = 167577
$ws.Range("C7").Value = 5.65
$ws.Range("C8").Value = 64.88
